# Applies the "ignore a bunch of stuff" edit: append four new job-action
# rows (98-101) to Sheet1 of the job-search tracking workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing data row (97) down into the four
# new rows so the new date cells pick up the same date-number style (s="3")
# as every other row in column A, and so the new rows otherwise look like a
# natural continuation of the table.
$ws.Range("A97:G97").Copy($ws.Range("A98:G101"))

# Row 98 - LatentView Analytics
$ws.Range("A98").Value = 45703
$ws.Range("B98").Value = "LatentView Analytics"
$ws.Range("C98").Value = "Senior Data Scientist"
$ws.Range("D98").Value = "Market analyytics, but seattle, easy"
$ws.Range("E98").ClearContents()
$ws.Range("F98").Value = "https://www.linkedin.com/jobs/view/4139907275/?refId=ByteString(length%3D16%2Cbytes%3D150983f7...8f6a262d)&trackingId=JQaQxdwWX0sVJpgbRDE9zw%3D%3D"
$ws.Range("G98").ClearContents()

# Row 99 - Tata Consultancy Services
$ws.Range("A99").Value = 45703
$ws.Range("B99").Value = "Tata Consultancy Services"
$ws.Range("C99").Value = "Data Scientist"
$ws.Range("D99").Value = "they want SQL… consulting"
$ws.Range("E99").ClearContents()
$ws.Range("F99").Value = "https://www.linkedin.com/jobs/view/4149289471/?refId=ByteString(length%3D16%2Cbytes%3Da77c4f9b...24e62942)&trackingId=fZD94CNeR4eK%2F1ctw%2B1IBg%3D%3D"
$ws.Range("G99").ClearContents()

# Row 100 - GLX ANALYTIX
$ws.Range("A100").Value = 45703
$ws.Range("C100").Value = "Senior Data Scientist "
$ws.Range("B100").Value = "GLX ANALYTIX"
$ws.Range("D100").Value = "Denmark, personalized medicine"
$ws.Range("E100").ClearContents()
$ws.Range("F100").Value = "https://www.linkedin.com/jobs/view/4150250844/?refId=ByteString(length%3D16%2Cbytes%3Db8b66b1f...f79b0827)&trackingId=IJSdoryuUWCZ8UQNmm2Csw%3D%3D"
$ws.Range("G100").ClearContents()

# Row 101 - SureCost
$ws.Range("A101").Value = 45703
$ws.Range("B101").Value = "SureCost"
$ws.Range("C101").Value = "Senior Data Scientist"
$ws.Range("D101").Value = "St. Petersburg, Fl,  pharmacy inventory, more of MLOPS?"
$ws.Range("E101").ClearContents()
$ws.Range("F101").Value = "https://www.linkedin.com/jobs/view/4152258208/?refId=I3ESEstkRIyQ0GNGtFx%2FYQ%3D%3D&trackingId=jyT38KXWQ46%2FKZ%2BxwTqcfQ%3D%3D"
$ws.Range("G101").ClearContents()

# Update the frozen-pane view + selection to match where the author ended up
# after adding the new rows.
$ws.Activate()
[void]$ws.Range("F101").Select()
